$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.66000000000057
$ws.Range("H2").Value = 0.0000433450846029082
$ws.Range("I2").Value = 0.0000433450846029082
$ws.Range("L2").Value = 48.67481866116682
$ws.Range("M2").Value = "[26.28859745358899, 71.06103986874466]"
$ws.Range("N2").Value = 0.00007025854592312974
$ws.Range("O2").Value = 0.00007025854592312974
$ws.Range("P2").Value = 1.616395018964117
$ws.Range("Q2").Value = "[1.0000264903318854, 2.232763547596349]"
$ws.Range("R2").Value = 0.000003579216438565069
$ws.Range("S2").Value = 0.000003579216438565069
$ws.Range("T2").Value = 64.18974815894937
$ws.Range("U2").Value = "[49.937325236121225, 78.44217108177752]"
$ws.Range("V2").Value = 0.00000000001005928673691869
$ws.Range("W2").Value = 0.00000000001005928673691869
$ws.Range("X2").Value = 19.05877877877921
$ws.Range("Y2").Value = 16.54158158158195
$ws.Range("Z2").Value = 21.57597597597646
$ws.Range("F3").Value = 25.66000000000057
$ws.Range("H3").Value = 0.000008549548495495607
$ws.Range("I3").Value = 0.000008549548495495607
$ws.Range("L3").Value = 54.00960579105124
$ws.Range("M3").Value = "[27.907523110053702, 80.11168847204877]"
$ws.Range("N3").Value = 0.0001379562071093243
$ws.Range("O3").Value = 0.0001379562071093243
$ws.Range("P3").Value = 1.025184389459733
$ws.Range("Q3").Value = "[0.5471843060306547, 1.5031844728888117]"
$ws.Range("R3").Value = 0.0000850422465028533
$ws.Range("S3").Value = 0.0000850422465028533
$ws.Range("T3").Value = 63.03400959557922
$ws.Range("U3").Value = "[49.1128793830965, 76.95513980806194]"
$ws.Range("V3").Value = 0.000000000008591793942969161
$ws.Range("W3").Value = 0.000000000008591793942969161
$ws.Range("X3").Value = 21.47323323323371
$ws.Range("Y3").Value = 19.52112112112155
$ws.Range("Z3").Value = 23.42534534534587
$ws.Range("B4").Value = 0
$ws.Range("F4").Value = 25.66000000000057
$ws.Range("H4").Value = 0.001168612257352875
$ws.Range("I4").Value = 0.001168612257352875
$ws.Range("L4").Value = 48.7723248806798
$ws.Range("M4").Value = "[15.618000278754835, 81.92664948260477]"
$ws.Range("N4").Value = 0.004855768529598681
$ws.Range("O4").Value = 0.004855768529598681
$ws.Range("P4").Value = 0.6226580034141929
$ws.Range("Q4").Value = "[-0.04402632347373192, 1.2893423303021176]"
$ws.Range("R4").Value = 0.0664384078697613
$ws.Range("S4").Value = 0.0664384078697613
$ws.Range("T4").Value = 66.14079226520606
$ws.Range("U4").Value = "[49.08626308430509, 83.19532144610703]"
$ws.Range("V4").Value = 0.0000000006501410521053685
$ws.Range("W4").Value = 0.0000000006501410521053685
$ws.Range("X4").Value = 23.11711711711763
$ws.Range("Y4").Value = 20.39443443443489
$ws.Range("Z4").Value = 25.83979979980038
$ws.Range("F5").Value = 25.66000000000057
$ws.Range("H5").Value = 0.01444873099678223
$ws.Range("I5").Value = 0.01444873099678223
$ws.Range("L5").Value = 33.51829432283844
$ws.Range("M5").Value = "[7.757109277709347, 59.279479367967525]"
$ws.Range("N5").Value = 0.01192807554714137
$ws.Range("O5").Value = 0.01192807554714137
$ws.Range("P5").Value = -0.0503157982556921
$ws.Range("Q5").Value = "[-1.0880791372793475, 0.9874475407679633]"
$ws.Range("R5").Value = 0.9226411107125156
$ws.Range("S5").Value = 0.9226411107125156
$ws.Range("T5").Value = 59.46103208626666
$ws.Range("U5").Value = "[44.25209829757221, 74.6699658749611]"
$ws.Range("V5").Value = 0.0000000005254958690414924
$ws.Range("W5").Value = 0.0000000005254958690414924
$ws.Range("X5").Value = 0.205485485485486
$ws.Range("Y5").Value = -4.032652652652747
$ws.Range("Z5").Value = 4.443623623623719
$ws.Range("F6").Value = 25.66000000000057
$ws.Range("H6").Value = 0.00001961626787272408
$ws.Range("I6").Value = 0.00001961626787272408
$ws.Range("L6").Value = 53.49175211627237
$ws.Range("M6").Value = "[28.759805889925588, 78.22369834261916]"
$ws.Range("N6").Value = 0.00007566060806452057
$ws.Range("O6").Value = 0.00007566060806452057
$ws.Range("P6").Value = -0.2012631930227693
$ws.Range("Q6").Value = "[-0.7421580242714629, 0.33963163822592435]"
$ws.Range("R6").Value = 0.4574958517167738
$ws.Range("S6").Value = 0.4574958517167738
$ws.Range("T6").Value = 61.5182266431347
$ws.Range("U6").Value = "[47.36651856108428, 75.66993472518513]"
$ws.Range("V6").Value = 0.00000000002815880861817277
$ws.Range("W6").Value = 0.00000000002815880861817277
$ws.Range("X6").Value = 0.821941941941958
$ws.Range("Y6").Value = -1.387027027027063
$ws.Range("Z6").Value = 3.030910910910979
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 25.66000000000057
$ws.Range("H7").Value = 0.000468459572188884
$ws.Range("I7").Value = 0.000468459572188884
$ws.Range("L7").Value = 37.41273523372926
$ws.Range("M7").Value = "[14.566464674473096, 60.259005792985434]"
$ws.Range("N7").Value = 0.00190630969592398
$ws.Range("O7").Value = 0.00190630969592398
$ws.Range("P7").Value = -0.3522105877898465
$ws.Range("Q7").Value = "[-1.0314738642416934, 0.32705268866200043]"
$ws.Range("R7").Value = 0.3019014184586286
$ws.Range("S7").Value = 0.3019014184586286
$ws.Range("T7").Value = 56.80067163485395
$ws.Range("U7").Value = "[44.45133259437655, 69.15001067533134]"
$ws.Range("V7").Value = 0.000000000005394351632048711
$ws.Range("W7").Value = 0.000000000005394351632048711
$ws.Range("X7").Value = 1.43839839839843
$ws.Range("Y7").Value = -1.335655655655685
$ws.Range("Z7").Value = 4.212452452452546
$ws.Range("F8").Value = 25.66000000000057
$ws.Range("H8").Value = 0.002331777581335404
$ws.Range("I8").Value = 0.002331777581335404
$ws.Range("L8").Value = 32.26910646711382
$ws.Range("M8").Value = "[7.479599914520186, 57.058613019707444]"
$ws.Range("N8").Value = 0.01189105925537914
$ws.Range("O8").Value = 0.01189105925537914
$ws.Range("P8").Value = -0.8302106712189241
$ws.Range("Q8").Value = "[-1.484316048542925, -0.17610529389492324]"
$ws.Range("R8").Value = 0.0140227501178205
$ws.Range("S8").Value = 0.0140227501178205
$ws.Range("T8").Value = 67.8212786812839
$ws.Range("U8").Value = "[55.20870577714203, 80.43385158542577]"
$ws.Range("V8").Value = 0.0000000000000404121180963557
$ws.Range("W8").Value = 0.0000000000000404121180963557
$ws.Range("X8").Value = 3.390510510510584
$ws.Range("Y8").Value = 0.7191991991992137
$ws.Range("Z8").Value = 6.061821821821955
$ws.Range("F9").Value = 23.34000000000021
$ws.Range("H9").Value = 0.000002925909526330095
$ws.Range("I9").Value = 0.000002925909526330095
$ws.Range("L9").Value = 44.75533326655705
$ws.Range("M9").Value = "[25.73233050386198, 63.778336029252124]"
$ws.Range("N9").Value = 0.00002184983593522816
$ws.Range("O9").Value = 0.00002184983593522816
$ws.Range("P9").Value = -1.295631805084079
$ws.Range("Q9").Value = "[-1.7736318885131563, -0.8176317216550011]"
$ws.Range("R9").Value = 0.000001966704909817096
$ws.Range("S9").Value = 0.000001966704909817096
$ws.Range("T9").Value = 64.4511484599544
$ws.Range("U9").Value = "[53.69304749628589, 75.20924942362291]"
$ws.Range("V9").Value = 0.000000000000001110223024625157
$ws.Range("W9").Value = 0.000000000000001110223024625157
$ws.Range("X9").Value = 4.812852852852895
$ws.Range("Y9").Value = 3.037237237237264
$ws.Range("Z9").Value = 6.588468468468527
$ws.Range("F10").Value = 23.34000000000021
$ws.Range("H10").Value = 0.02505987474205251
$ws.Range("I10").Value = 0.02505987474205251
$ws.Range("L10").Value = 30.28558771093495
$ws.Range("M10").Value = "[5.365466606790001, 55.205708815079895]"
$ws.Range("N10").Value = 0.01834038835927121
$ws.Range("O10").Value = 0.01834038835927121
$ws.Range("P10").Value = -1.534631846798618
$ws.Range("Q10").Value = "[-2.6667373075516974, -0.4025263860455386]"
$ws.Range("R10").Value = 0.009002987177638255
$ws.Range("S10").Value = 0.009002987177638255
$ws.Range("T10").Value = 68.64044169097015
$ws.Range("U10").Value = "[53.60692234040002, 83.67396104154028]"
$ws.Range("V10").Value = 0.000000000006713074540698472
$ws.Range("W10").Value = 0.000000000006713074540698472
$ws.Range("X10").Value = 5.700660660660713
$ws.Range("Y10").Value = 1.495255255255267
$ws.Range("Z10").Value = 9.90606606606616
$ws.Range("B11").Value = 1
$ws.Range("F11").Value = 23.34000000000021
$ws.Range("H11").Value = 0.003920855856629712
$ws.Range("I11").Value = 0.003920855856629712
$ws.Range("L11").Value = 36.41069255197863
$ws.Range("M11").Value = "[8.459564406860977, 64.36182069709629]"
$ws.Range("N11").Value = 0.01183459955197752
$ws.Range("O11").Value = 0.01183459955197752
$ws.Range("P11").Value = -1.006315965113847
$ws.Range("Q11").Value = "[-1.7610529389492324, -0.2515789912784623]"
$ws.Range("R11").Value = 0.01010638472671777
$ws.Range("S11").Value = 0.01010638472671777
$ws.Range("T11").Value = 52.03123779313636
$ws.Range("U11").Value = "[37.47866871706737, 66.58380686920535]"
$ws.Range("V11").Value = 0.000000005130520452922838
$ws.Range("W11").Value = 0.000000005130520452922838
$ws.Range("X11").Value = 3.738138138138169
$ws.Range("Y11").Value = 0.9345345345345422
$ws.Range("Z11").Value = 6.541741741741795
